# 03_MASTER_TAG_LIST.xlsx — "added tool and CLI"
#
# 1) Header row (A1:J1) loses its bold font + thin border + centered/top
#    alignment — revert back to the plain default cell style.
# 2) Many rows in column D ("target_name_description") get their
#    classification text normalised to "Unclassified" (title case) instead
#    of the previous "UNCLASSIFIED" / "Process Value" / "Control Valve"
#    placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) strip the header formatting (bold font / thin border / alignment) ---
$ws.Range("A1:J1").ClearFormats()

# --- 2) normalise column D classification text on the affected rows ---
$rows = @() + (16..19) + (25..29) + (59..63) + (67..69) + (72..78) + (92..93) + 130 + (143..145) + 149 + 152 + 162 + 164 + 170 + 172 + (174..181) + (207..208) + (213..214) + (221..224) + (226..240) + (248..252) + (255..258) + 338 + (342..350) + (352..361) + (363..364) + (366..405)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "Unclassified"
}
